$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 11 was missing a value in column E; fill it in with 8, matching the
# other plain numeric totals already present in that row (A11:H11).
$ws.Range("E11").Value = 8

# I11 already holds =SUM(A11:H11), so it recalculates from 57 to 65
# automatically once E11 is populated.
